$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.271983861923218
$ws.Range("B1").Value = 1.859593868255615
$ws.Range("C1").Value = 4.280735492706299
$ws.Range("D1").Value = 0.8278422355651855
$ws.Range("E1").Value = 0.794742226600647
